$d = $word.ActiveDocument

# --- Paragraph 1: "Test" -> bold "Iplc_sim_trap_address:" with double line spacing ---
$d.Content.Find.Execute("Test", $true, $false, $false, $false, $false, $true, 1, $false, `
    "Iplc_sim_trap_address:", 2) | Out-Null

$p1 = $d.Paragraphs(1)
$p1.Range.Font.Bold = $true
$p1.LineSpacingRule = 2

# --- Paragraph 2: insert the description run ahead of the existing bookmark, double spaced ---
$p2 = $d.Paragraphs(2)
$p2.Range.InsertBefore("This function deals with checking if a given address is in the cache. It takes into account the given associativity, and looks through the cache data structure. It will update the counter for a hit or a miss. After looking through the appropriate entries for the address it will call the appropriate function to deal with a hit or a miss.")
$p2.LineSpacingRule = 2
